$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.757.46'
$ws.Range("E2").Value = '  +0.28%  '

$ws.Range("D3").Value = '2.205.67'
$ws.Range("E3").Value = '  -1.23%  '

$ws.Range("E4").Value = '  +0.13%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '229.78'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.53%  '

$ws.Range("E6").Value = '  -0.76%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '60.45'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -1.96%  '

$ws.Range("E8").Value = '  +0.04%  '

$ws.Range("E9").Value = '  -0.16%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '56.96'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -3.25%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0889'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +1.28%  '

$ws.Range("E12").Value = '  -0.62%  '

$ws.Range("D13").Value = '2.533.36'
$ws.Range("E13").Value = '  -0.80%  '

$ws.Range("E14").Value = '  -1.70%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '22.10'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +0.34%  '

$ws.Range("E16").Value = '  -0.52%  '

$ws.Range("E17").Value = '  -0.08%  '

$ws.Range("D18").Value = '2.213.53'
$ws.Range("E18").Value = '  -0.91%  '

$ws.Range("D19").Value = '41.646.72'
$ws.Range("E19").Value = '  +0.43%  '

$ws.Range("D20").Value = '0.0₃0905'
$ws.Range("E20").Value = '  +0.26%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '71.90'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -1.85%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.03'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +0.37%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '241.48'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -2.33%  '

$ws.Range("E24").Value = '  -0.23%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.34'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -0.80%  '

$ws.Range("E26").Value = '  -4.85%  '

$ws.Range("E27").Value = '  +0.26%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '168.40'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -0.34%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.139'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -2.51%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.44'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +0.87%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '19.71'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -1.86%  '

$ws.Range("E32").Value = '  -7.39%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.120'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -1.24%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.98'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +0.44%  '

$ws.Range("E35").Value = '  -0.34%  '

$ws.Range("E36").Value = '  +3.14%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.55'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -5.59%  '

$ws.Range("E38").Value = '  -6.21%  '

$ws.Range("E39").Value = '  -1.92%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.000245'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +4.00%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.999'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +0.28%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0241'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +1.76%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.71'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +2.02%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0954'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -0.86%  '

$ws.Range("E45").Value = '  +1.25%  '

$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '96.44'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -3.41%  '

$ws.Range("B47").Value = 'FTXToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.36'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -9.88%  '

$ws.Range("D48").Value = '1.461.69'
$ws.Range("E48").Value = '  -1.81%  '

$ws.Range("E49").Value = '  -0.66%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '16.06'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -4.31%  '

$ws.Range("E51").Value = '  -1.05%  '
